# repair CPF/RG verification text ("{{RG}" -> "{{RG}}") and tidy up the
# "Sao Carlos, {{DATA}}" signature line on the certificate body (slide 1,
# shape "Rectangle 5").

$p  = $ppt.ActivePresentation
$s  = $p.Slides.Item(1)
$sh = $s.Shapes.Item(2)
$tr = $sh.TextFrame.TextRange

# The shape auto-fits its height to the text ("<a:spAutoFit/>"); remember the
# original height so it can be restored once the edits below are done, since
# the wording change should not resize the placeholder box.
$origHeight = $sh.Height

# --- Fix 1: "{{RG} " -> "{{RG}} " (missing closing brace before CPF) ------
$full = $tr.Text
$idxRG = $full.IndexOf("{{RG}")

# Leave the leading "{{" run untouched; rewrite the trailing "} " run (the
# last two characters of the "{{RG} " run) so it becomes "}} ", adding the
# missing closing curly brace.
$tail = $tr.Characters($idxRG + 5, 2)
$tail.Text = "}} "

# Re-assert the "RG" run text so it stands on its own run between the two
# brace runs, matching {{ | RG | }}<space>.
$mid = $tr.Characters($idxRG + 3, 2)
$mid.Text = "RG"

# --- Fix 2: merge "Sao Carlos" + ", " into a single run ------------------
$full2 = $tr.Text
$idxCarlos = $full2.IndexOf("Carlos") - 4
$idxData = $full2.IndexOf("{{DATA}}")
$cidade = $tr.Characters($idxCarlos + 1, $idxData - $idxCarlos)
$cidade.Text = $cidade.Text

# --- Restore the original (auto-fit) box height ---------------------------
$sh.Height = $origHeight
